$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column B, shifting dbExcel/WebExcel columns to C/D
$ws.Columns.Item(2).Insert()

# New header for the inserted "StatQuery" column
$ws.Range("B1").Value = "StatQuery"

# New stat query text for row 2, matching the wrap-text style used by A2
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma, NOS'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# New column width for inserted column (matches column A width)
$ws.Columns.Item(2).ColumnWidth = 75

# Update selection
$ws.Range("A2").Select()

$wb.Save()
